# Update automàtic: dades i banners [2026-02-19 19:20]
# Applies the scraped-data refresh to sheet "Dades_Meteo" (ActiveSheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-19 19:18:25"
$ws.Range("I2").Value = "2.7 mm"
$ws.Range("E3").Value = "2026-02-19 19:18:27"
$ws.Range("I3").Value = "4.1 mm"
$ws.Range("E4").Value = "2026-02-19 19:18:30"
$ws.Range("H4").Value = "'57%"
$ws.Range("J4").Value = "1009.6 hPa"
$ws.Range("O4").Value = "11.8 °C"
$ws.Range("E5").Value = "2026-02-19 19:18:32"
$ws.Range("I5").Value = "7.3 mm"
$ws.Range("E6").Value = "2026-02-19 19:18:34"
$ws.Range("J6").Value = "1009.7 hPa"
$ws.Range("E7").Value = "2026-02-19 19:18:37"
$ws.Range("J7").Value = "1010.6 hPa"
$ws.Range("E8").Value = "2026-02-19 19:18:39"
$ws.Range("J8").Value = "1010.3 hPa"
$ws.Range("O8").Value = "10.0 °C"
$ws.Range("E9").Value = "2026-02-19 19:18:42"
$ws.Range("O9").Value = "10.6 °C"
$ws.Range("E10").Value = "2026-02-19 19:18:44"
$ws.Range("H10").Value = "'69%"
$ws.Range("O10").Value = "10.8 °C"
$ws.Range("E11").Value = "2026-02-19 19:18:45"
$ws.Range("H11").Value = "'66%"
$ws.Range("O11").Value = "5.5 °C"
$ws.Range("E12").Value = "2026-02-19 19:18:46"
$ws.Range("E13").Value = "2026-02-19 19:18:47"
$ws.Range("J13").Value = "1010.9 hPa"
$ws.Range("O13").Value = "4.3 °C"
$ws.Range("E14").Value = "2026-02-19 19:18:49"
$ws.Range("E15").Value = "2026-02-19 19:18:50"
$ws.Range("O15").Value = "10.1 °C"
$ws.Range("E16").Value = "2026-02-19 19:18:51"
$ws.Range("I16").Value = "8.3 mm"
$ws.Range("E17").Value = "2026-02-19 19:18:52"
$ws.Range("H17").Value = "'82%"
$ws.Range("E18").Value = "2026-02-19 19:18:53"
$ws.Range("J18").Value = "1009.9 hPa"
$ws.Range("E19").Value = "2026-02-19 19:18:54"
$ws.Range("H19").Value = "'77%"
$ws.Range("O19").Value = "5.4 °C"
$ws.Range("E20").Value = "2026-02-19 19:18:55"
$ws.Range("E21").Value = "2026-02-19 19:18:56"
$ws.Range("H21").Value = "'62%"
$ws.Range("J21").Value = "1010.8 hPa"
$ws.Range("E22").Value = "2026-02-19 19:18:59"
$ws.Range("E23").Value = "2026-02-19 19:19:01"
$ws.Range("G23").Value = "212 cm"
$ws.Range("I23").Value = "8.8 mm"
$ws.Range("E24").Value = "2026-02-19 19:19:04"
$ws.Range("J24").Value = "1014.4 hPa"
$ws.Range("E25").Value = "2026-02-19 19:19:06"
$ws.Range("I25").Value = "4.6 mm"
$ws.Range("E26").Value = "2026-02-19 19:19:09"
$ws.Range("J26").Value = "1009.6 hPa"
$ws.Range("E27").Value = "2026-02-19 19:19:12"
$ws.Range("O27").Value = "-3.7 °C"
$ws.Range("E28").Value = "2026-02-19 19:19:14"
$ws.Range("H28").Value = "'67%"
$ws.Range("J28").Value = "1009.5 hPa"
$ws.Range("L28").Value = "27.4 km/h - 273º 18:59 TU"
$ws.Range("O28").Value = "9.2 °C"
$ws.Range("E29").Value = "2026-02-19 19:19:17"
$ws.Range("E30").Value = "2026-02-19 19:19:19"
$ws.Range("J30").Value = "1009.7 hPa"
$ws.Range("E31").Value = "2026-02-19 19:19:22"
$ws.Range("J31").Value = "1009.1 hPa"
$ws.Range("E32").Value = "2026-02-19 19:19:24"
$ws.Range("E33").Value = "2026-02-19 19:19:27"
$ws.Range("L33").Value = "55.4 km/h - 310º 18:33 TU"
$ws.Range("O33").Value = "3.6 °C"
$ws.Range("E34").Value = "2026-02-19 19:19:30"
$ws.Range("E35").Value = "2026-02-19 19:19:32"
$ws.Range("H35").Value = "'67%"
$ws.Range("J35").Value = "1015.9 hPa"
$ws.Range("E36").Value = "2026-02-19 19:19:35"
$ws.Range("J36").Value = "1010.0 hPa"
$ws.Range("E37").Value = "2026-02-19 19:19:37"
$ws.Range("H37").Value = "'70%"
$ws.Range("J37").Value = "1011.0 hPa"
$ws.Range("O37").Value = "5.9 °C"
$ws.Range("E38").Value = "2026-02-19 19:19:40"
$ws.Range("H38").Value = "'57%"
$ws.Range("E39").Value = "2026-02-19 19:19:42"
$ws.Range("H39").Value = "'73%"
$ws.Range("I39").Value = "4.5 mm"
$ws.Range("E40").Value = "2026-02-19 19:19:44"
$ws.Range("H40").Value = "'73%"
$ws.Range("J40").Value = "1012.0 hPa"
$ws.Range("O40").Value = "6.3 °C"
$ws.Range("E41").Value = "2026-02-19 19:19:47"
$ws.Range("J41").Value = "1012.6 hPa"
$ws.Range("E42").Value = "2026-02-19 19:19:49"
$ws.Range("O42").Value = "11.4 °C"
$ws.Range("E43").Value = "2026-02-19 19:19:52"
$ws.Range("E44").Value = "2026-02-19 19:19:54"
$ws.Range("I44").Value = "8.5 mm"
$ws.Range("E45").Value = "2026-02-19 19:19:57"
$ws.Range("J45").Value = "1015.0 hPa"
$ws.Range("E46").Value = "2026-02-19 19:20:00"
$ws.Range("J46").Value = "1015.3 hPa"
